# Fruta / hortaliza, semanal
# Insert a new weekly record at row 41 (shifting existing rows 41-60 down to 42-61)
# for "Vega Monumental Concepción" / Ciruela / Angeleno / Segunda.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 41; existing rows 41..60 shift down to 42..61.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new data record.
$ws.Range("A41").Value = 11
$ws.Range("B41").Value = "Vega Monumental Concepción"
$ws.Range("C41").Value = "Bíobío"
$ws.Range("D41").Value = 44636
$ws.Range("E41").Value = 8
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100103
$ws.Range("H41").Value = "Frutos de hueso (carozo)"
$ws.Range("I41").Value = 100103002
$ws.Range("J41").Value = "Ciruela"
$ws.Range("K41").Value = "Angeleno"
$ws.Range("L41").Value = "Segunda"
$ws.Range("M41").Value = 220
$ws.Range("N41").Value = 6500
$ws.Range("O41").Value = 7000
$ws.Range("P41").Value = 6727
$ws.Range("Q41").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R41").Value = "Provincia de Curicó"
$ws.Range("S41").Value = 374
$ws.Range("T41").Value = 18
